$wb = $excel.ActiveWorkbook

# --- Sheet "2o Parcial": update rows 9-12, columns E:K ---
$ws2 = $wb.Worksheets.Item("2o Parcial")

# Row 9
$ws2.Cells.Item(9, 5).Value  = 11     # E9  Aprobados
$ws2.Cells.Item(9, 6).Value  = 0      # F9  Reprobados
$ws2.Cells.Item(9, 7).Value  = 100    # G9  por_aprobados
$ws2.Cells.Item(9, 8).Value  = 0      # H9  por_reprobados
$ws2.Cells.Item(9, 9).Value  = 8.6    # I9  Promedio
$ws2.Cells.Item(9, 10).Value = 0      # J9  Blancos
$ws2.Cells.Item(9, 11).Value = 0      # K9  por_blancos

# Row 10
$ws2.Cells.Item(10, 5).Value  = 24
$ws2.Cells.Item(10, 6).Value  = 1
$ws2.Cells.Item(10, 7).Value  = 96
$ws2.Cells.Item(10, 8).Value  = 4
$ws2.Cells.Item(10, 9).Value  = 9.4
$ws2.Cells.Item(10, 10).Value = 0
$ws2.Cells.Item(10, 11).Value = 0

# Row 11
$ws2.Cells.Item(11, 5).Value  = 35
$ws2.Cells.Item(11, 6).Value  = 1
$ws2.Cells.Item(11, 7).Value  = 97.2
$ws2.Cells.Item(11, 8).Value  = 2.8
$ws2.Cells.Item(11, 9).Value  = 9
$ws2.Cells.Item(11, 10).Value = 0
$ws2.Cells.Item(11, 11).Value = 0

# Row 12
$ws2.Cells.Item(12, 5).Value  = 89
$ws2.Cells.Item(12, 6).Value  = 5
$ws2.Cells.Item(12, 7).Value  = 94.7
$ws2.Cells.Item(12, 8).Value  = 5.3
$ws2.Cells.Item(12, 9).Value  = 8.699999999999999
$ws2.Cells.Item(12, 10).Value = 0
$ws2.Cells.Item(12, 11).Value = 0

# --- Sheet "Final": update Promedio column (I) for rows 9-12 ---
$ws3 = $wb.Worksheets.Item("Final")

$ws3.Cells.Item(9, 9).Value  = 8.6
$ws3.Cells.Item(10, 9).Value = 8.5
$ws3.Cells.Item(11, 9).Value = 8.6
$ws3.Cells.Item(12, 9).Value = 8.5
